$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 16 (shape 3) switches from the custom "Table_0" style to
#    PowerPoint's built-in "Medium Style 2 - Accent 1" table style.
# ---------------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$tbl = $s16.Shapes.Item(3).Table
$tbl.ApplyStyle("{7624B20A-C123-4E4B-AB3C-3B01365B60B9}")

# ---------------------------------------------------------------------------
# 2) The deck's theme colour scheme (the one used by every slide via the
#    slide master) swaps from the "Integral" palette to the stock "Office"
#    palette. RGB values below are the standard Office theme colours,
#    encoded as the BGR integers the PowerPoint object model expects.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$scheme = $s1.ThemeColorScheme

$scheme.Colors(1).RGB  = 0          # Dark 1    -> 000000
$scheme.Colors(2).RGB  = 16777215   # Light 1   -> FFFFFF
$scheme.Colors(3).RGB  = 6968388    # Dark 2    -> 44546A
$scheme.Colors(4).RGB  = 15132391   # Light 2   -> E7E6E6
$scheme.Colors(5).RGB  = 13998939   # Accent 1  -> 5B9BD5
$scheme.Colors(6).RGB  = 3243501    # Accent 2  -> ED7D31
$scheme.Colors(7).RGB  = 10855845   # Accent 3  -> A5A5A5
$scheme.Colors(8).RGB  = 49407      # Accent 4  -> FFC000
$scheme.Colors(9).RGB  = 12874308   # Accent 5  -> 4472C4
$scheme.Colors(10).RGB = 4697456    # Accent 6  -> 70AD47
$scheme.Colors(11).RGB = 12673797   # Hyperlink -> 0563C1
$scheme.Colors(12).RGB = 7491477    # Followed Hyperlink -> 954F72
